# Apply the "Update CDA Logical model for ST.r2b" edit.
#
# Changes:
#  1. Metadata sheet: bump Version + Date strings.
#  2. Metadata sheet: insert a new "Jurisdiction" row (empty value) right
#     after "Contact" and before "Description" (formatting copied from the
#     row above so it matches the other data rows).
#  3. Elements sheet: add the II-1 invariant text to the Constraint(s)
#     column (AJ) for the AssociatedEntity.typeId row (row 5).

$wb = $excel.ActiveWorkbook

# ---- 1 & 2: Metadata sheet ------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Update Version (row 3) and Date (row 8) values.
$meta.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$meta.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" right before "Description" (currently
# row 11), pushing Description and everything below it down by one.
$meta.Rows.Item(11).Insert()

# Copy the formatting from the row above (Contact, row 10) so the new row
# matches the rest of the table's style.
$meta.Range("A10:B10").Copy()
$meta.Range("A11:B11").PasteSpecial(-4122)

$meta.Range("A11").Value = "Jurisdiction"
$meta.Range("B11").Value = ""

# ---- 3: Elements sheet ----------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 5 is AssociatedEntity.typeId (Type(s) = II). Add the II-1 constraint.
$elements.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}
"
